$d = $word.ActiveDocument

# "possible fix for mals2-34": the merge-field placeholder {d.PreviousMonth}
# should reference the date the results were actually reported on, not the
# previous month, so the field name changes from PreviousMonth to
# ReportedOnDate.
$d.Content.Find.Execute("PreviousMonth", $true, $false, $false, $false,
                         $false, $true, 1, $false, "ReportedOnDate", 2)
